$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 3; $r -le 46; $r++) {
    $ws.Cells.Item($r, 2).Value = 123456
}

$ws.Range("G10").Select()
